$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking values
# (e.g. "533.23") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range('D2:D51')
$priceRange.NumberFormat = '@'

$ws.Range('D2').Value = '58.701.01'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.512.19'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '533.23'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '136.38'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('D9').Value = '2.514.04'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  +1.02%  '
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('D13').Value = '0.345'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = '2.937.13'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').Value = '23.02'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = '58.714.38'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '2.505.94'
$ws.Range('E18').Value = '  -1.66%  '
$ws.Range('D19').Value = '10.99'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '4.23'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '323.89'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = '5.85'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').Value = '63.48'
$ws.Range('E24').Value = '  +2.53%  '
$ws.Range('D25').Value = '0.416'
$ws.Range('E25').Value = '  -0.91%  '
$ws.Range('D26').Value = '0.164'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('D28').Value = '7.55'
$ws.Range('E28').Value = '  -3.62%  '
$ws.Range('D29').Value = '6.71'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').Value = '0.0₃0768'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('D31').Value = '1.76'
$ws.Range('E31').Value = '  -2.05%  '
$ws.Range('D32').Value = '166.26'
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').Value = '1.15'
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').Value = '1.38'
$ws.Range('E35').Value = '  -4.74%  '
$ws.Range('D36').Value = '18.43'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').Value = '4.08'
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').Value = '36.65'
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('D40').Value = '0.815'
$ws.Range('E40').Value = '  +1.13%  '
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').Value = '5.22'
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('D43').Value = '277.35'
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '0.599'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').Value = '10.84'
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('D47').Value = '126.41'
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('D48').Value = '0.0923'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('D49').Value = '0.0510'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '0.0220'
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('D51').Value = '17.32'
$ws.Range('E51').Value = '  -1.38%  '

# Restore the default (Normal) cell style now that the text values are set,
# so the cells end up with no explicit style index, same as the source file.
$priceRange.Style = 'Normal'
